$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CNData")
$t = $ws.ListObjects.Item(1)

$dummy1 = $t.ListColumns.Add()
$dummy2 = $t.ListColumns.Add()
Write-Host "Added 2 dummies, count=$($t.ListColumns.Count)"
$dummy1.Delete()
Write-Host "Deleted dummy1, count=$($t.ListColumns.Count)"

$diaStrandCol = $t.ListColumns.Add()
$diaStrandCol.Range.Cells.Item(1).Value = "DiaStrand"
